$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing-row tweaks (rows 92-94): DataType column ---
$ws.Cells.Item(92, 9).Value = "CBA,Cells"
$ws.Cells.Item(93, 9).Value = "Cells"
$ws.Cells.Item(94, 9).Value = "CBA,Cells"

# ===================================================================
# The shared-string table in the target file is appended to in the
# exact order the cells below were authored; reproduce that order
# precisely so the resulting xl/sharedStrings.xml matches byte-for-byte
# in content (even though the surrounding rows are filled in out of
# that order further down).
# ===================================================================

# (1)(2)(3) Name column for the three new Kinetics rows (99-101)
$ws.Cells.Item(99, 3).Value = "B16OVAIFNgPulsedKineticsPart1"
$ws.Cells.Item(100, 3).Value = "B16OVAIFNgPulsedKineticsPart2"
$ws.Cells.Item(101, 3).Value = "B16OVAIFNgPulsedKineticsPart3"

# (4) Row 97 comment gets extended
$ws.Cells.Item(97, 18).Value = "last three timepoint supernantants dried out due to low incubator humidity, calibration's undiluted sample was weird for some reason; ended up discarding the sample for 15 total calibration solutions"

# (5)(6) Row 103 Name + ExperimentType
$ws.Cells.Item(103, 3).Value = "322and412_CBA_Dilutions"
$ws.Cells.Item(103, 5).Value = "DilutionCorrection"

# (7)(8) Row 102 Name + Full Name (literal, not a formula)
$ws.Cells.Item(102, 3).Value = "B16OVAIFNgPulsedConfluency"
$ws.Cells.Item(102, 7).Value = "20190423-B16OVAIFNgPulsedConfluency_B16_Timeseries_1"

# (9) Row 103 comment
$ws.Cells.Item(103, 18).Value = "1:5 Dilution not enough for N4 1uM IL2; will have to redo"

# (10) Row 104 Name
$ws.Cells.Item(104, 3).Value = "404and308"

# ===================================================================
# Fill in the remaining cells for rows 97 and 99-104 (re-using existing
# shared strings / plain numbers, so ordering no longer matters).
# ===================================================================

# Row 97: new DataType / CBA dilution columns
$ws.Cells.Item(97, 9).Value = "CBA,Cells"
$ws.Cells.Item(97, 12).Value = 15
$ws.Cells.Item(97, 13).Value = 1

# Row 99: B16OVAIFNgPulsedKineticsPart1
$ws.Cells.Item(99, 1).Value = 98
$ws.Cells.Item(99, 2).Value = 20190416
$ws.Cells.Item(99, 4).Value = "B16"
$ws.Cells.Item(99, 4).Font.Color = 0
$ws.Cells.Item(99, 5).Value = "Timeseries"
$ws.Cells.Item(99, 6).Value = 1
$ws.Cells.Item(99, 8).Value = "Emanuel"
$ws.Cells.Item(99, 9).Value = "Cells"
$ws.Cells.Item(99, 10).Value = 3
$ws.Cells.Item(99, 11).Value = 8
$ws.Cells.Item(99, 17).Value = "hand"

# Row 100: B16OVAIFNgPulsedKineticsPart2
$ws.Cells.Item(100, 1).Value = 99
$ws.Cells.Item(100, 2).Value = 20190416
$ws.Cells.Item(100, 4).Value = "B16"
$ws.Cells.Item(100, 4).Font.Color = 0
$ws.Cells.Item(100, 5).Value = "Timeseries"
$ws.Cells.Item(100, 6).Value = 1
$ws.Cells.Item(100, 8).Value = "Emanuel"
$ws.Cells.Item(100, 9).Value = "Cells"
$ws.Cells.Item(100, 10).Value = 3
$ws.Cells.Item(100, 11).Value = 48
$ws.Cells.Item(100, 17).Value = "hand"

# Row 101: B16OVAIFNgPulsedKineticsPart3
$ws.Cells.Item(101, 1).Value = 100
$ws.Cells.Item(101, 2).Value = 20190416
$ws.Cells.Item(101, 4).Value = "B16"
$ws.Cells.Item(101, 4).Font.Color = 0
$ws.Cells.Item(101, 5).Value = "Timeseries"
$ws.Cells.Item(101, 6).Value = 1
$ws.Cells.Item(101, 8).Value = "Emanuel"
$ws.Cells.Item(101, 9).Value = "Cells"
$ws.Cells.Item(101, 10).Value = 5
$ws.Cells.Item(101, 11).Value = 48
$ws.Cells.Item(101, 17).Value = "hand"

# Row 102: B16OVAIFNgPulsedConfluency
$ws.Cells.Item(102, 1).Value = 101
$ws.Cells.Item(102, 2).Value = 20190423
$ws.Cells.Item(102, 4).Value = "B16"
$ws.Cells.Item(102, 5).Value = "Timeseries"
$ws.Cells.Item(102, 6).Value = 1
$ws.Cells.Item(102, 8).Value = "Emanuel"
$ws.Cells.Item(102, 9).Value = "Cells"
$ws.Cells.Item(102, 10).Value = 1
$ws.Cells.Item(102, 11).Value = 96
$ws.Cells.Item(102, 17).Value = "hand"

# Row 103: 322and412_CBA_Dilutions
$ws.Cells.Item(103, 1).Value = 102
$ws.Cells.Item(103, 2).Value = 20190426
$ws.Cells.Item(103, 4).Value = "OT1"
$ws.Cells.Item(103, 4).Font.Color = 0
$ws.Cells.Item(103, 6).Value = 1
$ws.Cells.Item(103, 8).Value = "Sooraj"
$ws.Cells.Item(103, 9).Value = "CBA"
$ws.Cells.Item(103, 10).Value = 36
$ws.Cells.Item(103, 11).Value = 5
$ws.Cells.Item(103, 12).Value = 16
$ws.Cells.Item(103, 13).Value = 1
$ws.Cells.Item(103, 17).Value = "hand"

# Row 104: 404and308
$ws.Cells.Item(104, 1).Value = 103
$ws.Cells.Item(104, 2).Value = 20190429
$ws.Cells.Item(104, 4).Value = "OT1"
$ws.Cells.Item(104, 4).Font.Color = 0
$ws.Cells.Item(104, 5).Value = "DilutionCorrection"
$ws.Cells.Item(104, 6).Value = 2
$ws.Cells.Item(104, 8).Value = "Sooraj"
$ws.Cells.Item(104, 9).Value = "CBA"
$ws.Cells.Item(104, 10).Value = 36
$ws.Cells.Item(104, 11).Value = 16
$ws.Cells.Item(104, 12).Value = 16
$ws.Cells.Item(104, 13).Value = 1
$ws.Cells.Item(104, 17).Value = "hand"

# --- Full Name formula column (G): extend the shared formula down
# through the new rows (skipping row 102, which keeps a literal value) ---
$ws.Range("G99:G101").Formula = '=B99&"-"&C99&"_"&D99&"_"&E99&"_"&F99'
$ws.Range("G103:G104").Formula = '=B103&"-"&C103&"_"&D103&"_"&E103&"_"&F103'

# --- View: update selection to match the end-state cursor position ---
$ws.Range("G116").Select() | Out-Null

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
